$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New accuracy values for column B, keyed by row number (retrained model run)
$newValues = @{
    2 = 0.9375
    3 = 0.84375
    4 = 0.796875
    6 = 0.65625
    7 = 0.59375
    8 = 0.578125
    9 = 0.65625
    10 = 0.609375
    11 = 0.578125
    12 = 0.625
    13 = 0.625
    14 = 0.65625
    15 = 0.609375
    16 = 0.59375
    17 = 0.625
    18 = 0.640625
    19 = 0.609375
    20 = 0.671875
    21 = 0.59375
    22 = 0.5625
    23 = 0.4375
    25 = 0.453125
    26 = 0.484375
    27 = 0.515625
    28 = 0.46875
    29 = 0.5
    30 = 0.4375
    31 = 0.453125
    33 = 0.484375
    35 = 0.484375
    36 = 0.484375
    37 = 0.484375
    38 = 0.46875
    39 = 0.46875
    40 = 0.46875
    41 = 0.46875
    42 = 0.46875
    43 = 0.46875
    44 = 0.46875
    69 = 0.453125
    70 = 0.453125
    71 = 0.453125
    72 = 0.453125
    73 = 0.453125
    74 = 0.453125
    75 = 0.453125
    76 = 0.453125
    77 = 0.453125
    78 = 0.453125
    79 = 0.453125
    80 = 0.453125
    81 = 0.453125
    82 = 0.453125
    83 = 0.453125
    84 = 0.453125
    85 = 0.453125
    86 = 0.453125
    87 = 0.453125
    88 = 0.46875
    89 = 0.46875
    90 = 0.46875
    91 = 0.46875
    92 = 0.46875
    93 = 0.46875
    94 = 0.46875
    95 = 0.46875
    96 = 0.46875
    97 = 0.46875
    98 = 0.46875
    99 = 0.46875
    100 = 0.46875
    101 = 0.46875
    102 = 0.46875
    103 = 0.484375
    104 = 0.53125
    105 = 0.390625
    106 = 0.375
    107 = 0.4375
    109 = 0.546875
    110 = 0.34375
    111 = 0.4375
    112 = 0.3125
    113 = 0.484375
    114 = 0.515625
    115 = 0.4375
    116 = 0.53125
    117 = 0.453125
    118 = 0.4426229508196721
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 2).Value = $newValues[$row]
}

# The notebook object repr in column A picks up the new runs memory address
$newLabel = "<__main__.DisplayOutputs object at 0x7f11b010d4f0>"
for ($row = 102; $row -le 118; $row++) {
    $ws.Cells.Item($row, 1).Value = $newLabel
}

# Restore the selection to the data range, matching the post-edit workbook state
[void]$ws.Range("A2:B118").Select()

